$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$values = @{
    "A2" = 45875
    "B2" = 81.34999999999999
    "C2" = 80.31999999999999
    "D2" = 80
    "E2" = 80.40000000000001
    "F2" = 85.72
    "G2" = 97.34999999999999
    "H2" = 94
    "I2" = 104.41
    "J2" = 97.34999999999999
    "K2" = 76.95
    "L2" = 51
    "M2" = 27.67
    "N2" = 28.4
    "O2" = 27.2
    "P2" = 27.2
    "Q2" = 27.2
    "R2" = 28.8
    "S2" = 52
    "T2" = 80.31999999999999
    "U2" = 106.22
    "V2" = 114.64
    "W2" = 170
    "X2" = 150.01
    "Y2" = 115.43
    "Z2" = 78.5
    "AA2" = "20h-24h"
    "AB2" = 137.52
    "AC2" = "20h-22h"
    "AD2" = 142.32
    "AE2" = "22h-24h"
    "AF2" = 132.72
    "AG2" = "9h-17h"
}

foreach ($addr in $values.Keys) {
    $ws.Range($addr).Value = $values[$addr]
}

$wb.Save()
